$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the other
# header cells (bold, bordered) by copying formatting from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate the time_taken values for each data row (F2:F12).
$timestamps = @(
    "2021-10-05 10:51:49.600312",
    "2021-10-05 10:51:49.600325",
    "2021-10-05 10:51:49.600329",
    "2021-10-05 10:51:49.600332",
    "2021-10-05 10:51:49.600336",
    "2021-10-05 10:51:49.600339",
    "2021-10-05 10:51:49.600342",
    "2021-10-05 10:51:49.600345",
    "2021-10-05 10:51:49.600349",
    "2021-10-05 10:51:49.600352",
    "2021-10-05 10:51:49.600355"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
